$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - Wins, Losses, Ties in columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting from an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Data rows 2-43: Wins=79, Losses=83, Ties=0
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 79   # AD
    $ws.Cells.Item($r, 31).Value = 83   # AE
    $ws.Cells.Item($r, 32).Value = 0    # AF
}
